# Reformat the `questions = [...]` Python-literal payload in A2 into a
# pretty-printed JSON-ish block, move it up to A1 (replacing the old
# placeholder "0" value), drop the now-empty second row, and strip the
# bold/centered/bordered formatting that used to sit on A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You have a String word in your Java code, and you want to obtain a count of the number of characters in it. Which version of the code accomplishes this?",
        "ques_type": 2,
        "options": [
            "word.length()",
            "word.length",
            "word.size",
            "word.size()"
        ],
        "score": "word.length()"
    },
    {
        "title": "You want to create an object of a class without providing any additional information to the constructor.  Which type of constructor do you need to use?",
        "ques_type": 2,
        "options": [
            "no-detail constructor",
            "default constructor",
            "no-argument constructor",
            "static constructor"
        ],
        "score": "no-argument constructor"
    },
    {
        "title": "You want to declare an integer variable var1 in a class such that it will be visible only inside that class. Which version of the code accomplishes this?",
        "ques_type": 2,
        "options": [
            "local int var1",
            "final int var1",
            "package int var1",
            "private int var1"
        ],
        "score": "private int var1"
    },
    {
        "title": "True or false: One method can contain multiple generic arguments.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "True"
    }
]
'@

# A1 currently holds a bordered/bold/centered style (xfId 1); clear it back
# to the workbook default before writing the new text so the cell picks up
# no formatting, matching the plain (un-styled) target cell.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# The multi-line text would otherwise leave a custom row height behind;
# auto-fit row 1 back to the sheet's normal height.
$ws.Rows.Item(1).AutoFit()

# Row 2 (the old shared-string cell) is no longer needed now that its text
# lives in A1 - remove the whole row.
$ws.Range("A2").EntireRow.Delete()
